$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "211.02") must be forced to text
# so Excel does not auto-convert them to numbers, matching the original
# inlineStr string storage used throughout this sheet.
$textForceCells = @(
    "D5", "D9", "D10", "D11", "D14", "D16", "D19", "D21", "D22", "D24", "D25", "D27", "D29", "D33", "D35", "D39", "D42", "D43", "D44", "D46", "D47", "D48"
)
foreach ($cellAddr in $textForceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Apply all updated values from the cryptos refresh
$ws.Range("D2").Value = "26.640.84"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "1.591.77"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "211.02"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "0.0616"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "1.814.72"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "1.590.51"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "64.63"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "26.635.52"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "208.61"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "4.25"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "146.84"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "7.27"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").Value = "0.671"
$ws.Range("E33").Value = "  +23.70%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.322.46"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "2.91"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("E36").Value = "  -4.25%  "
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D39").Value = "0.831"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").Value = "0.789"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").Value = "63.13"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "1.727.52"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "90.02"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "0.833"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -0.48%  "

# Restore default (Normal) style on the text-forced cells so no stray
# number-format styling is introduced by this edit.
foreach ($cellAddr in $textForceCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
